$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = 'magapoke_2026-02-04'

$ws.Cells.Item(1,1).Value = 'rank'
$ws.Cells.Item(1,2).Value = 'title'

$titles = @(
    '黄昏町プリズナーズ',
    'ハンドレッドノート－アグリーダック－',
    '黒月のイェルクナハト',
    'ドリーム☆ジャンボ☆ガール',
    'K-9~警視庁公安部公安第9課異能対策係~',
    '【爆アド】生まれた直後から最強悪霊と脳内バトルしてたら魔力量が測定可能域を超えてました〜悪憑の子の謙虚な覇道〜',
    'アイドラトリィ',
    'せいぶつ部の田辺くん',
    '篝家の８兄弟',
    'ハードワーカー中田',
    'あの島の海音荘',
    'ともだちづくり',
    'ナキナギ',
    'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！',
    '追放されなかった男　～二度目の人生は土下座から始まりました～',
    '普通の本はありません！',
    'ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜',
    '邪目さんは邪神です',
    '屋根の下のアルテミス',
    'ゼロとヒャク',
    '春くらり',
    '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～',
    '白鳥運子は31画',
    'MYS',
    '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～',
    'その青春',
    '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～',
    '君が監督！',
    '鳴るさんだぁ',
    '明智ナンバーワン',
    'ハプスブルク家の華麗なる受難',
    'じゅーくぼっくす',
    'ナマイキ旭ちゃんをわからせたい',
    '歪みの虜',
    '平成転生',
    'JK Biker',
    '宇曽田みのりの代用料理',
    '夜鐘のキト',
    'GURU',
    '永久のユウグレ',
    '人生逆転ダンジョン',
    '花子狩り',
    '眠れる森のレガ',
    '白銀のキュイジーヌ～明治外交官の料理人～',
    'きゃわるり方程式',
    '〈小市民〉 春期限定いちごタルト事件',
    'イエティ、とある日々'
)

for ($i = 0; $i -lt $titles.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $titles[$i]
}

Write-Host "Done. Sheets: $($wb.Worksheets.Count)"
